$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes ---------------------------------------------------
# Insert a new "Phase Name" column before the existing "Description" column
# (this pushes the old "Description"/"Duration(days)" columns one to the
# right, which is exactly what the new layout needs).
$ws.Columns("E").Insert()

# --- Copy the existing body formatting into the new column and the two
#     extra Gantt-phase rows that are being added below the table ----------
$ws.Range("F4").Copy() | Out-Null
$ws.Range("E4:E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C11:G11").Copy() | Out-Null
$ws.Range("C12:G13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Header row ------------------------------------------------------------
$ws.Range("C3").Value2 = "Start Date"
$ws.Range("D3").Value2 = "End Date"
$ws.Range("E3").Value2 = "Phase Name"
$ws.Range("F3").Value2 = "Description"
$ws.Range("G3").Value2 = "Duration(days)"

# --- Data rows ---------------------------------------------------------
$data = @(
  @{ Row=4;  C=44472; D=44607; E="Research ";                 F="Research Phase" },
  @{ Row=5;  C=44500; D=44520; E="Desing ";                    F="Design GUI" },
  @{ Row=6;  C=44514; D=44531; E="GUI Skeleton";                F="Code GUI skeleton to start adding the functions to the program" },
  @{ Row=7;  C=44532; D=44555; E="First Function";              F="Code the usages times function" },
  @{ Row=8;  C=44544; D=44557; E="Individual Testing / Fix";    F="Test the first function and if have errors correct it." },
  @{ Row=9;  C=44557; D=44582; E="Second and Third Function";   F="Code the schedule function and program lock function" },
  @{ Row=10; C=44564; D=44588; E="Functions Testing /Fix";      F="Test the 2 last functions, and if the have errors correct it." },
  @{ Row=11; C=44585; D=44613; E="Finalize the GUI ";           F="Add the Themes to the GUI, and stylize it." },
  @{ Row=12; C=44557; D=44625; E="User Test ";                  F="Test the program in its entirety and give it as a beta program to some selected users with the purpose of testing it" },
  @{ Row=13; C=44627; D=44635; E="Fix";                         F="Fix all the possible errors that the program may have" }
)

foreach ($d in $data) {
  $r = $d.Row
  $ws.Range("C$r").Value2 = $d.C
  $ws.Range("D$r").Value2 = $d.D
  $ws.Range("E$r").Value2 = $d.E
  $ws.Range("F$r").Value2 = $d.F
}

# Duration(days) column: a single formula, then filled down (creates the
# shared-formula group Excel itself would produce).
$ws.Range("G4").Formula = "=DAYS(D4,C4)"
$ws.Range("G5:G13").Formula = "=DAYS(D5,C5)"

# --- Wrap text on the long phase-name cells that need it -------------------
$ws.Range("E9").WrapText = $true
$ws.Range("E10").WrapText = $true
$ws.Range("E12").WrapText = $true

# --- Column widths -------------------------------------------------------
$ws.Columns("E").ColumnWidth = 25.5703125
$ws.Columns("F").ColumnWidth = 103.85546875
$ws.Columns("G").ColumnWidth = 16.5703125

# --- Chart: extend the data ranges to the new rows/columns -----------------
$chart = $ws.ChartObjects(1).Chart
$chart.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$C`$3,Sheet1!`$E`$4:`$E`$13,Sheet1!`$C`$4:`$C`$13,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(Sheet1!`$G`$3,Sheet1!`$E`$4:`$E`$13,Sheet1!`$G`$4:`$G`$13,2)"

# Move the chart down two rows to stay below the now-larger data table.
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + ($ws.Rows("12").RowHeight + $ws.Rows("13").RowHeight)

# --- Selection -------------------------------------------------------------
$ws.Range("H23").Select() | Out-Null
